# Ticker appears with results: add a new data row (row 3) to the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (number format / style) from the existing row 2 cells that
# carry a non-default style, then overwrite with the new row's values.
$ws.Range("A2").Copy($ws.Range("A3"))
$ws.Range("S2").Copy($ws.Range("S3"))

$ws.Range("A3").Value = 42632.882164351853
$ws.Range("B3").Value = -7
$ws.Range("C3").Value = "Neutral"
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = "Random"
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 1.77
$ws.Range("S3").Value = 0.1132
$ws.Range("T3").Value = -4.05
$ws.Range("U3").Value = 5.85
$ws.Range("V3").Value = "N/A"
$ws.Range("W3").Value = 0
